# Weekly fruit/vegetable data update:
# Insert 3 new rows (new price-report week, date 45267) at the top of the
# "Femacal de La Calera" data block, pushing the existing rows down by 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("665:667").Insert()

# New row 665 - Especial
$ws.Range("A665").Value = 3
$ws.Range("B665").Value = "Femacal de La Calera"
$ws.Range("C665").Value = "Coquimbo"
$ws.Range("D665").Value = 45267
$ws.Range("E665").Value = 5
$ws.Range("F665").Value = "Fruta"
$ws.Range("G665").Value = 100101
$ws.Range("H665").Value = "Berries"
$ws.Range("I665").Value = 100112025
$ws.Range("J665").Value = "Frutilla"
$ws.Range("K665").Value = "Sin especificar"
$ws.Range("L665").Value = "Especial"
$ws.Range("M665").Value = 60
$ws.Range("N665").Value = 10000
$ws.Range("O665").Value = 10000
$ws.Range("P665").Value = 10000
$ws.Range("Q665").Value = "`$/bandeja 7 kilos"
$ws.Range("R665").Value = "Provincia de Melipilla"
$ws.Range("S665").Value = 1429
$ws.Range("T665").Value = 7

# New row 666 - Primera
$ws.Range("A666").Value = 3
$ws.Range("B666").Value = "Femacal de La Calera"
$ws.Range("C666").Value = "Coquimbo"
$ws.Range("D666").Value = 45267
$ws.Range("E666").Value = 5
$ws.Range("F666").Value = "Fruta"
$ws.Range("G666").Value = 100101
$ws.Range("H666").Value = "Berries"
$ws.Range("I666").Value = 100112025
$ws.Range("J666").Value = "Frutilla"
$ws.Range("K666").Value = "Sin especificar"
$ws.Range("L666").Value = "Primera"
$ws.Range("M666").Value = 67
$ws.Range("N666").Value = 8000
$ws.Range("O666").Value = 8000
$ws.Range("P666").Value = 8000
$ws.Range("Q666").Value = "`$/bandeja 7 kilos"
$ws.Range("R666").Value = "Provincia de Melipilla"
$ws.Range("S666").Value = 1143
$ws.Range("T666").Value = 7

# New row 667 - Segunda
$ws.Range("A667").Value = 3
$ws.Range("B667").Value = "Femacal de La Calera"
$ws.Range("C667").Value = "Coquimbo"
$ws.Range("D667").Value = 45267
$ws.Range("E667").Value = 5
$ws.Range("F667").Value = "Fruta"
$ws.Range("G667").Value = 100101
$ws.Range("H667").Value = "Berries"
$ws.Range("I667").Value = 100112025
$ws.Range("J667").Value = "Frutilla"
$ws.Range("K667").Value = "Sin especificar"
$ws.Range("L667").Value = "Segunda"
$ws.Range("M667").Value = 60
$ws.Range("N667").Value = 6000
$ws.Range("O667").Value = 6000
$ws.Range("P667").Value = 6000
$ws.Range("Q667").Value = "`$/bandeja 7 kilos"
$ws.Range("R667").Value = "Provincia de Melipilla"
$ws.Range("S667").Value = 857
$ws.Range("T667").Value = 7
